$d = $word.ActiveDocument
$rng = $d.Content
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00123909" w:rsidRDefault="00B658B3" w:rsidP="00B658B3"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">At first, I was very intimidated by this assignment. The scale and complexity of that one page is daunting. </w:t></w:r><w:r><w:t xml:space="preserve">It includes a fixed menu bar, several promotion items with multiple images and links, as well as a detailed footer. All with the ability to respond to screen size changes in a dramatic way. </w:t></w:r><w:r><w:t xml:space="preserve">However, after many hours I think my version of the website is very close (aesthetically) to the real version. Essentially none of the links work, since I didn’t see that as an important aspect of the assignment, but almost all the of the intricate details of the real webpage can be seen on my own. For </w:t></w:r><w:r w:rsidR="00764702"><w:t>example,</w:t></w:r><w:r><w:t xml:space="preserve"> hovering over the menu bar changes the different buttons to a grey color and I believe all the correct links will appear with an underline when hovered over. </w:t></w:r><w:r w:rsidR="00764702"><w:t>Also,</w:t></w:r><w:r><w:t xml:space="preserve"> the responsiveness of my version resembles the real page.</w:t></w:r></w:p><w:p w:rsidR="00B658B3" w:rsidRDefault="00B658B3" w:rsidP="00B658B3"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p w:rsidR="00B658B3" w:rsidRDefault="00B658B3" w:rsidP="00B658B3"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>The most complex aspect to implement w</w:t></w:r><w:r><w:t>as</w:t></w:r><w:r><w:t xml:space="preserve"> the formatting </w:t></w:r><w:r w:rsidR="002C751B"><w:t>of</w:t></w:r><w:r><w:t xml:space="preserve"> the various sections of the page. For </w:t></w:r><w:r w:rsidR="00764702"><w:t>example,</w:t></w:r><w:r><w:t xml:space="preserve"> the menu bar was very challenging since everything is supposed to be centered and responsive. What I noticed throughout the assignment was that </w:t></w:r><w:r w:rsidR="00764702"><w:t xml:space="preserve">when I first started I wasn’t very good at writing the proper HTML or CSS code. However, after </w:t></w:r><w:r w:rsidR="003452F9"><w:t>a few</w:t></w:r><w:r><w:t xml:space="preserve"> hours of struggling in the beginning I learned from my mistakes and began to really utilize the inspect tool within chrome. I was able to quickly diagnose how the real website was written and recreate it in my own </w:t></w:r><w:r><w:t>code</w:t></w:r><w:r><w:t>. The other aspect that I struggled with early on was proper class names. Early in the assignment I wasn’t writing very clean or efficient code and I’m sure I doubled up on some CSS code that could have easily been captured with better class naming structure. However, when I was working on the later stages of the page I realized that the naming and structure to the HTML could save a lot of work and I think I improved my efficiency as I worked my way through the assignment.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">The only known bugs or differences between my webpage and the real page is the links as I mentioned before. Some of the links </w:t></w:r><w:r><w:t xml:space="preserve">on the real website </w:t></w:r><w:r><w:t>also cover whole sections (including background photos) and I couldn’t figure out how to get the whole section to be a clickable link even though my &lt;a&gt; tag spanned the HTML</w:t></w:r><w:r><w:t xml:space="preserve"> content</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">I really enjoyed this assignment. I feel that I learned a ton from building this page. In my opinion </w:t></w:r><w:r><w:t>not</w:t></w:r><w:r><w:t xml:space="preserve"> everything in HTML or CSS is very intuitive until you actually work with it and learn how to play with the code (maybe that’s true of all programming</w:t></w:r><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">:) </w:t></w:r><w:r><w:t>)</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>. I thought this was a great assignment, lots of work but also lots of reward to see myself learning and building something</w:t></w:r><w:r><w:t xml:space="preserve"> so visual</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>Thank you,</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>Aaron Gilman</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
